$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.70643966666667
$ws.Range("H2").Value = 53.119319
$ws.Range("I2").Value = 0.4380235920947999
$ws.Range("J2").Value = 0.4380235920947999
$ws.Range("M2").Value = 46.29121633333333
$ws.Range("N2").Value = 138.873649
$ws.Range("O2").Value = 0.3133663986859022
$ws.Range("P2").Value = 0.3133663986859022
$ws.Range("Q2").Value = 819.6526291027812
$ws.Range("R2").Value = 7376.873661925031
$ws.Range("S2").Value = 0.13726187559421
$ws.Range("T2").Value = 0.13726187559421
$ws.Range("G3").Value = 17.70643966666667
$ws.Range("H3").Value = 53.119319
$ws.Range("I3").Value = 0.4380235920947999
$ws.Range("J3").Value = 0.4380235920947999
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("O3").Value = 0.3169204109998198
$ws.Range("P3").Value = 0.3169204109998198
$ws.Range("Q3").Value = 828.948633872859
$ws.Range("R3").Value = 7460.53770485573
$ws.Range("S3").Value = 0.1388186168343014
$ws.Range("T3").Value = 0.1388186168343014
$ws.Range("G4").Value = 17.70643966666667
$ws.Range("H4").Value = 53.119319
$ws.Range("I4").Value = 0.4380235920947999
$ws.Range("J4").Value = 0.4380235920947999
$ws.Range("M4").Value = 38.53544233333333
$ws.Range("N4").Value = 115.606327
$ws.Range("O4").Value = 0.2608640200510233
$ws.Range("P4").Value = 0.2608640200510233
$ws.Range("Q4").Value = 682.3254847034792
$ws.Range("R4").Value = 6140.929362331313
$ws.Range("S4").Value = 0.1142645951110391
$ws.Range("T4").Value = 0.1142645951110391
$ws.Range("G5").Value = 17.70643966666667
$ws.Range("H5").Value = 53.119319
$ws.Range("I5").Value = 0.4380235920947999
$ws.Range("J5").Value = 0.4380235920947999
$ws.Range("M5").Value = 16.07945366666667
$ws.Range("N5").Value = 48.238361
$ws.Range("O5").Value = 0.1088491702632547
$ws.Range("P5").Value = 0.1088491702632547
$ws.Range("Q5").Value = 284.7098762217955
$ws.Range("R5").Value = 2562.388885996159
$ws.Range("S5").Value = 0.04767850455524931
$ws.Range("T5").Value = 0.04767850455524931
$ws.Range("G6").Value = 1.617245333333334
$ws.Range("H6").Value = 4.851736000000001
$ws.Range("I6").Value = 0.04000756919748267
$ws.Range("J6").Value = 0.04000756919748267
$ws.Range("M6").Value = 46.29121633333333
$ws.Range("N6").Value = 138.873649
$ws.Range("O6").Value = 0.3133663986859022
$ws.Range("P6").Value = 0.3133663986859022
$ws.Range("Q6").Value = 74.86425358940713
$ws.Range("R6").Value = 673.7782823046641
$ws.Range("S6").Value = 0.01253702787959217
$ws.Range("T6").Value = 0.01253702787959217
$ws.Range("G7").Value = 1.617245333333334
$ws.Range("H7").Value = 4.851736000000001
$ws.Range("I7").Value = 0.04000756919748267
$ws.Range("J7").Value = 0.04000756919748267
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("O7").Value = 0.3169204109998198
$ws.Range("P7").Value = 0.3169204109998198
$ws.Range("S7").Value = 0.01267921527316994
$ws.Range("T7").Value = 0.01267921527316994
$ws.Range("G8").Value = 1.617245333333334
$ws.Range("H8").Value = 4.851736000000001
$ws.Range("I8").Value = 0.04000756919748267
$ws.Range("J8").Value = 0.04000756919748267
$ws.Range("M8").Value = 38.53544233333333
$ws.Range("N8").Value = 115.606327
$ws.Range("O8").Value = 0.2608640200510233
$ws.Range("P8").Value = 0.2608640200510233
$ws.Range("Q8").Value = 62.32126428151911
$ws.Range("R8").Value = 560.8913785336721
$ws.Range("S8").Value = 0.01043653533332482
$ws.Range("T8").Value = 0.01043653533332482
$ws.Range("G9").Value = 1.617245333333334
$ws.Range("H9").Value = 4.851736000000001
$ws.Range("I9").Value = 0.04000756919748267
$ws.Range("J9").Value = 0.04000756919748267
$ws.Range("M9").Value = 16.07945366666667
$ws.Range("N9").Value = 48.238361
$ws.Range("O9").Value = 0.1088491702632547
$ws.Range("P9").Value = 0.1088491702632547
$ws.Range("Q9").Value = 26.00442140496623
$ws.Range("R9").Value = 234.039792644696
$ws.Range("S9").Value = 0.004354790711395736
$ws.Range("T9").Value = 0.004354790711395736
$ws.Range("G10").Value = 21.099799
$ws.Range("H10").Value = 63.299397
$ws.Range("I10").Value = 0.5219688387077175
$ws.Range("J10").Value = 0.5219688387077175
$ws.Range("M10").Value = 46.29121633333333
$ws.Range("N10").Value = 138.873649
$ws.Range("O10").Value = 0.3133663986859022
$ws.Range("P10").Value = 0.3133663986859022
$ws.Range("Q10").Value = 976.7353600988504
$ws.Range("R10").Value = 8790.618240889653
$ws.Range("S10").Value = 0.1635674952121
$ws.Range("T10").Value = 0.1635674952121
$ws.Range("G11").Value = 21.099799
$ws.Range("H11").Value = 63.299397
$ws.Range("I11").Value = 0.5219688387077175
$ws.Range("J11").Value = 0.5219688387077175
$ws.Range("M11").Value = 46.81622333333333
$ws.Range("O11").Value = 0.3169204109998198
$ws.Range("P11").Value = 0.3169204109998198
$ws.Range("Q11").Value = 987.8129022724434
$ws.Range("R11").Value = 8890.316120451989
$ws.Range("S11").Value = 0.1654225788923485
$ws.Range("T11").Value = 0.1654225788923485
$ws.Range("G12").Value = 21.099799
$ws.Range("H12").Value = 63.299397
$ws.Range("I12").Value = 0.5219688387077175
$ws.Range("J12").Value = 0.5219688387077175
$ws.Range("M12").Value = 38.53544233333333
$ws.Range("N12").Value = 115.606327
$ws.Range("O12").Value = 0.2608640200510233
$ws.Range("P12").Value = 0.2608640200510233
$ws.Range("Q12").Value = 813.0900876094242
$ws.Range("R12").Value = 7317.810788484819
$ws.Range("S12").Value = 0.1361628896066594
$ws.Range("T12").Value = 0.1361628896066594
$ws.Range("G13").Value = 21.099799
$ws.Range("H13").Value = 63.299397
$ws.Range("I13").Value = 0.5219688387077175
$ws.Range("J13").Value = 0.5219688387077175
$ws.Range("M13").Value = 16.07945366666667
$ws.Range("N13").Value = 48.238361
$ws.Range("O13").Value = 0.1088491702632547
$ws.Range("P13").Value = 0.1088491702632547
$ws.Range("Q13").Value = 339.2732403964797
$ws.Range("R13").Value = 3053.459163568317
$ws.Range("S13").Value = 0.05681587499660969
$ws.Range("T13").Value = 0.05681587499660969
